$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug reports")

# Set cell values (reuses the existing "Link to attachments" shared string)
$ws.Range("J7").Value = "Link to attachments"
$ws.Range("J8").Value = "Link to attachments"
$ws.Range("J9").Value = "Link to attachments"

# Add hyperlinks pointing at the BR6/BR7/BR8 attachment folders, following
# the same convention as the existing BR4/BR5 hyperlinks
$ws.Hyperlinks.Add($ws.Range("J7"), "https://github.com/Oleksandr-Mnk/Test-documentation/tree/main/Bug%20reports/Attachments%20to%20bug%20reports/BR6")
$ws.Hyperlinks.Add($ws.Range("J8"), "https://github.com/Oleksandr-Mnk/Test-documentation/tree/main/Bug%20reports/Attachments%20to%20bug%20reports/BR7")
$ws.Hyperlinks.Add($ws.Range("J9"), "https://github.com/Oleksandr-Mnk/Test-documentation/tree/main/Bug%20reports/Attachments%20to%20bug%20reports/BR8")

# Re-apply the same formatting (border/hyperlink font) used by the other
# cells in column J, since adding a hyperlink resets the cell style
$ws.Range("J6").Copy()
$ws.Range("J7:J9").PasteSpecial(-4122)

# Clear clipboard marquee / selection like Excel leaves after a copy/paste
$excel.CutCopyMode = $false
$ws.Range("I11").Select()
